$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.268.93"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "1.822.96"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'313.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "'0.4469"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.65%  "

$ws.Range("D8").Value = "'0.3764"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("D9").Value = "'0.07402"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").Value = "'0.8794"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.71%  "

$ws.Range("D11").Value = "'20.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").Value = "1.822.24"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D14").Value = "'5.427"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.89%  "

$ws.Range("D15").Value = "'92.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").Value = "'0.07059"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "'0.000008812"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").Value = "'15.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("D21").Value = "27.268.14"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("D22").Value = "'5.339"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'10.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").Value = "'1.961"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Value = "'150.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").Value = "'2.278"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.56%  "

$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("D28").Value = "'5.349"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "

$ws.Range("D29").Value = "'117.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").Value = "'0.08889"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").Value = "'0.7887"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.84%  "

$ws.Range("D32").Value = "'1.197"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").Value = "'4.576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.83%  "

$ws.Range("D34").Value = "'2.921"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "'0.9995"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").Value = "'1.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.51%  "

$ws.Range("D37").Value = "'0.01973"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("D38").Value = "'0.05272"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "

$ws.Range("D39").Value = "'7.293"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.41%  "

$ws.Range("D40").Value = "'0.5299"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.347"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.23%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.873"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("D43").Value = "'0.1701"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").Value = "'8.635"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.80%  "

$ws.Range("D45").Value = "'0.5053"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.18%  "

$ws.Range("D46").Value = "'10.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").Value = "'105.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("E48").Value = "  +1.19%  "

$ws.Range("D49").Value = "'0.9991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").Value = "'66.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.87%  "
